$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-LinkCell($ws, $cell, $url) {
    # Write the URL as the cell's text value, then paint the whole string
    # (in two adjacent character runs so it stays a shared-string rich-text
    # run rather than collapsing into a cell-level style) with the
    # underline + link-blue formatting used throughout the sheet, and
    # finally register the real hyperlink relationship.
    $cell.Value = $url
    $len = $url.Length
    $run1 = $cell.Characters(1, $len - 1)
    $run1.Font.Underline = 2
    $run1.Font.ColorIndex = 4
    $run2 = $cell.Characters($len, 1)
    $run2.Font.Underline = 2
    $run2.Font.ColorIndex = 4
    $ws.Hyperlinks.Add($cell, $url, "", "", $url)
}

# ------------------------------------------------------------------
# 1) Grow the sheet: rows 200-226 become new blank rows (matching the
#    existing blank-row look of rows 197-199), giving room for the six
#    new data rows (197-202) plus new trailing blank rows (203-226).
# ------------------------------------------------------------------
$ws.Range("A200:E226").RowHeight = 15
$ws.Range("A199:E199").Copy()
$ws.Range("A200:E226").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 2) Fill in the six new event rows (197-202).
# ------------------------------------------------------------------
$ws.Range("A197").Value = 45710
$ws.Range("B197").Value = "TECHNOALLIANZ SCHRANZ ONLY"
$ws.Range("C197").Value = "Elektroküche"
$ws.Range("D197").Value = "Köln"
Set-LinkCell $ws $ws.Range("E197") "https://www.instagram.com/reel/DGIfNeQgYYt/?igsh=aG50bTdhYmVkZW9k"

$ws.Range("A198").Value = 45709
$ws.Range("B198").Value = "THREE SINS CLUB HARDTECHNO"
$ws.Range("C198").Value = "Oma Doris"
$ws.Range("D198").Value = "Dortmund"
Set-LinkCell $ws $ws.Range("E198") "https://www.instagram.com/p/DF--9scqFLS/?igsh=ZGgweWxwbDJyMXJl"

$ws.Range("A199").Value = 45715
$ws.Range("B199").Value = "TURBO 120 MIN RAVE"
$ws.Range("C199").Value = "Oma Doris"
$ws.Range("D199").Value = "Dortmund"
Set-LinkCell $ws $ws.Range("E199") "https://www.instagram.com/reel/DGLoCJ3qENM/?igsh=MXZpZThlenV3cm53ZQ=="

$ws.Range("A200").Value = 45718
$ws.Range("B200").Value = "BEATS BASS CARNIVAL"
$ws.Range("C200").Value = "Odonien"
$ws.Range("D200").Value = "Köln"
Set-LinkCell $ws $ws.Range("E200") "https://www.instagram.com/reel/DFsGr5os7I8/?igsh=ZDl5dWMxZ2lycWZ2"

$ws.Range("A201").Value = 45724
$ws.Range("B201").Value = "TECHNOBLOCK"
$ws.Range("C201").Value = "Elektroküche"
$ws.Range("D201").Value = "Köln"
Set-LinkCell $ws $ws.Range("E201") "https://www.instagram.com/reel/DGDxCLCsRH6/?igsh=NHpyNG5pc2NoY25w"

$ws.Range("A202").Value = 45777
$ws.Range("B202").Value = "MYRAVE RAVE IN DEN MAI"
$ws.Range("C202").Value = "Essigfabrik & Elektroküche"
$ws.Range("D202").Value = "Köln"
Set-LinkCell $ws $ws.Range("E202") "https://www.instagram.com/reel/DGBa_TeA81I/?igsh=dWUydzdmYXlycGZh"

# ------------------------------------------------------------------
# 3) Re-apply the normal data-row formatting (border/fill/number format)
#    to the six rows we just populated, since writing .Value alone keeps
#    whatever formatting the cell already had (the blank-row look).
# ------------------------------------------------------------------
$ws.Range("A196:E196").Copy()
$ws.Range("A197:E202").PasteSpecial(-4122)
